$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) updates to column F ("想去人数")
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 672
$wsExhibit.Range("F9").Value  = 1056
$wsExhibit.Range("F15").Value = 523
$wsExhibit.Range("F16").Value = 527
$wsExhibit.Range("F22").Value = 437
$wsExhibit.Range("F25").Value = 1116
$wsExhibit.Range("F26").Value = 233
$wsExhibit.Range("F29").Value = 1175
$wsExhibit.Range("F30").Value = 430
$wsExhibit.Range("F32").Value = 3783
$wsExhibit.Range("F34").Value = 692

# Sheet "全部类型" (sheet4.xml) updates to column F ("想去人数")
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F14").Value = 672
$wsAll.Range("F15").Value = 672
$wsAll.Range("F17").Value = 1056
$wsAll.Range("F27").Value = 523
$wsAll.Range("F28").Value = 527
$wsAll.Range("F36").Value = 437
$wsAll.Range("F41").Value = 1116
$wsAll.Range("F42").Value = 233
$wsAll.Range("F47").Value = 1175
$wsAll.Range("F48").Value = 430
$wsAll.Range("F49").Value = 3783
$wsAll.Range("F52").Value = 692

$wb.Save()
